# ------------------------------------------------------------------
# "towards getting time series structure achieved"
#
# 1) proxy2: drop the d11Bsd (C) values on rows that are not full
#    replicate measurements (no d11B/B value), and reset the sheet
#    selection/active-cell state.
# 2) Add a new "proxy3" sheet holding the consolidated time series
#    (duplicate-age rows merged) and make it the active tab.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- proxy2 cleanup ---
$ws2 = $wb.Worksheets.Item("proxy2")
$clearCells = @("C3","C4","C5","C6","C7","C10","C11","C12","C13","C14","C15","C18","C19","C20","C21","C22","C23","C24","C28","C31","C32","C33","C34","C38","C39","C42","C43","C44","C45","C46")
foreach ($addr in $clearCells) {
    $ws2.Range($addr).ClearContents()
}

# --- add proxy3 as the new last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "proxy3"

# --- populate proxy3 with the consolidated time-series table ---
$ws3.Cells.Item(1,1).Value = "age"
$ws3.Cells.Item(1,2).Value = "d11B"
$ws3.Cells.Item(1,3).Value = "d11Bsd"
$ws3.Cells.Item(1,4).Value = "d18O"
$ws3.Cells.Item(1,5).Value = "MgCa"
$ws3.Cells.Item(1,6).Value = "species"

$ws3.Cells.Item(2,1).Value = 58.776
$ws3.Cells.Item(2,2).Value = 16.67
$ws3.Cells.Item(2,3).Value = 0.13
$ws3.Cells.Item(2,4).Value = -1.13
$ws3.Cells.Item(2,5).Value = 3.59
$ws3.Cells.Item(2,6).Value = "Grub"

$ws3.Cells.Item(3,1).Value = 58.621
$ws3.Cells.Item(3,4).Value = -1.13
$ws3.Cells.Item(3,5).Value = 3.5
$ws3.Cells.Item(3,6).Value = "Grub"

$ws3.Cells.Item(4,1).Value = 58.441
$ws3.Cells.Item(4,4).Value = -1.22
$ws3.Cells.Item(4,5).Value = 3.98
$ws3.Cells.Item(4,6).Value = "Grub"

$ws3.Cells.Item(5,1).Value = 58.286
$ws3.Cells.Item(5,4).Value = -1.07
$ws3.Cells.Item(5,5).Value = 3.82
$ws3.Cells.Item(5,6).Value = "Grub"

$ws3.Cells.Item(6,1).Value = 58.133
$ws3.Cells.Item(6,4).Value = -1.3
$ws3.Cells.Item(6,5).Value = 3.93
$ws3.Cells.Item(6,6).Value = "Grub"

$ws3.Cells.Item(7,1).Value = 58.102
$ws3.Cells.Item(7,4).Value = -1.59
$ws3.Cells.Item(7,5).Value = 4.37
$ws3.Cells.Item(7,6).Value = "Grub"

$ws3.Cells.Item(8,1).Value = 58.101
$ws3.Cells.Item(8,2).Value = 16.67
$ws3.Cells.Item(8,3).Value = 0.115
$ws3.Cells.Item(8,5).Value = 4.4
$ws3.Cells.Item(8,6).Value = "Grub"

$ws3.Cells.Item(9,1).Value = 58.099
$ws3.Cells.Item(9,2).Value = 16.22
$ws3.Cells.Item(9,3).Value = 0.14
$ws3.Cells.Item(9,5).Value = 4.16
$ws3.Cells.Item(9,6).Value = "Grub"

$ws3.Cells.Item(10,1).Value = 58.086
$ws3.Cells.Item(10,5).Value = 4.47
$ws3.Cells.Item(10,6).Value = "Grub"

$ws3.Cells.Item(11,1).Value = 57.987
$ws3.Cells.Item(11,2).Value = 16.47
$ws3.Cells.Item(11,3).Value = 0.13
$ws3.Cells.Item(11,4).Value = -1.65
$ws3.Cells.Item(11,5).Value = 3.88
$ws3.Cells.Item(11,6).Value = "Grub"

$ws3.Cells.Item(12,1).Value = 57.806
$ws3.Cells.Item(12,2).Value = 16.67
$ws3.Cells.Item(12,3).Value = 0.095
$ws3.Cells.Item(12,4).Value = -1.53
$ws3.Cells.Item(12,5).Value = 3.52
$ws3.Cells.Item(12,6).Value = "Grub"

$ws3.Cells.Item(13,1).Value = 57.619
$ws3.Cells.Item(13,4).Value = -1.6
$ws3.Cells.Item(13,5).Value = 3.89
$ws3.Cells.Item(13,6).Value = "Grub"

$ws3.Cells.Item(14,1).Value = 57.5
$ws3.Cells.Item(14,4).Value = -1.32
$ws3.Cells.Item(14,5).Value = 3.38
$ws3.Cells.Item(14,6).Value = "Grub"

$ws3.Cells.Item(15,1).Value = 57.338
$ws3.Cells.Item(15,4).Value = -1.62
$ws3.Cells.Item(15,5).Value = 3.67
$ws3.Cells.Item(15,6).Value = "Grub"

$ws3.Cells.Item(16,1).Value = 57.226
$ws3.Cells.Item(16,5).Value = 3.58
$ws3.Cells.Item(16,6).Value = "Grub"

$ws3.Cells.Item(17,1).Value = 57.074
$ws3.Cells.Item(17,4).Value = -1.42
$ws3.Cells.Item(17,5).Value = 3.77
$ws3.Cells.Item(17,6).Value = "Grub"

$ws3.Cells.Item(18,1).Value = 56.12
$ws3.Cells.Item(18,5).Value = 3.63
$ws3.Cells.Item(18,6).Value = "Grub"

$ws3.Cells.Item(19,1).Value = 55.965
$ws3.Cells.Item(19,2).Value = 15.76
$ws3.Cells.Item(19,3).Value = 0.215
$ws3.Cells.Item(19,4).Value = -1.35
$ws3.Cells.Item(19,5).Value = 3.9
$ws3.Cells.Item(19,6).Value = "Grub"

$ws3.Cells.Item(20,1).Value = 55.957
$ws3.Cells.Item(20,4).Value = -1.45
$ws3.Cells.Item(20,5).Value = 3.53
$ws3.Cells.Item(20,6).Value = "Grub"

$ws3.Cells.Item(21,1).Value = 55.956
$ws3.Cells.Item(21,2).Value = 15.46
$ws3.Cells.Item(21,3).Value = 0.11
$ws3.Cells.Item(21,4).Value = -1.55
$ws3.Cells.Item(21,5).Value = 3.54
$ws3.Cells.Item(21,6).Value = "Grub"

$ws3.Cells.Item(22,1).Value = 55.95
$ws3.Cells.Item(22,2).Value = 15.73
$ws3.Cells.Item(22,3).Value = 0.195
$ws3.Cells.Item(22,4).Value = -1.53
$ws3.Cells.Item(22,5).Value = 3.33
$ws3.Cells.Item(22,6).Value = "Grub"

$ws3.Cells.Item(23,1).Value = 55.946
$ws3.Cells.Item(23,4).Value = -1.62
$ws3.Cells.Item(23,5).Value = 3.35
$ws3.Cells.Item(23,6).Value = "Grub"

$ws3.Cells.Item(24,1).Value = 55.934
$ws3.Cells.Item(24,4).Value = -2.0
$ws3.Cells.Item(24,5).Value = 4.69
$ws3.Cells.Item(24,6).Value = "Grub"

$ws3.Cells.Item(25,1).Value = 55.932
$ws3.Cells.Item(25,4).Value = -1.99
$ws3.Cells.Item(25,5).Value = 5.04
$ws3.Cells.Item(25,6).Value = "Grub"

$ws3.Cells.Item(26,1).Value = 55.932
$ws3.Cells.Item(26,4).Value = -1.99
$ws3.Cells.Item(26,5).Value = 4.69
$ws3.Cells.Item(26,6).Value = "Grub"

$ws3.Cells.Item(27,1).Value = 55.927
$ws3.Cells.Item(27,2).Value = 14.82
$ws3.Cells.Item(27,3).Value = 0.15
$ws3.Cells.Item(27,4).Value = -2.0
$ws3.Cells.Item(27,5).Value = 5.37
$ws3.Cells.Item(27,6).Value = "Grub"

$ws3.Cells.Item(28,1).Value = 55.914
$ws3.Cells.Item(28,2).Value = 15.11
$ws3.Cells.Item(28,3).Value = 0.14
$ws3.Cells.Item(28,4).Value = -1.89
$ws3.Cells.Item(28,5).Value = 5.09
$ws3.Cells.Item(28,6).Value = "Grub"

$ws3.Cells.Item(29,1).Value = 55.901
$ws3.Cells.Item(29,2).Value = 14.45
$ws3.Cells.Item(29,3).Value = 0.225
$ws3.Cells.Item(29,4).Value = -1.81
$ws3.Cells.Item(29,5).Value = 4.83
$ws3.Cells.Item(29,6).Value = "Grub"

$ws3.Cells.Item(30,1).Value = 55.901
$ws3.Cells.Item(30,4).Value = -1.81
$ws3.Cells.Item(30,5).Value = 4.83
$ws3.Cells.Item(30,6).Value = "Grub"

$ws3.Cells.Item(31,1).Value = 55.888
$ws3.Cells.Item(31,4).Value = -1.75
$ws3.Cells.Item(31,5).Value = 5.01
$ws3.Cells.Item(31,6).Value = "Grub"

$ws3.Cells.Item(32,1).Value = 55.885
$ws3.Cells.Item(32,2).Value = 14.9
$ws3.Cells.Item(32,3).Value = 0.15
$ws3.Cells.Item(32,4).Value = -1.84
$ws3.Cells.Item(32,5).Value = 4.41
$ws3.Cells.Item(32,6).Value = "Grub"

$ws3.Cells.Item(33,1).Value = 55.872
$ws3.Cells.Item(33,2).Value = 15.09
$ws3.Cells.Item(33,3).Value = 0.12
$ws3.Cells.Item(33,4).Value = -1.57
$ws3.Cells.Item(33,5).Value = 4.79
$ws3.Cells.Item(33,6).Value = "Grub"

$ws3.Cells.Item(34,1).Value = 55.846
$ws3.Cells.Item(34,4).Value = -1.51
$ws3.Cells.Item(34,5).Value = 4.08
$ws3.Cells.Item(34,6).Value = "Grub"

$ws3.Cells.Item(35,1).Value = 55.815
$ws3.Cells.Item(35,4).Value = -1.49
$ws3.Cells.Item(35,5).Value = 3.83
$ws3.Cells.Item(35,6).Value = "Grub"

$ws3.Cells.Item(36,1).Value = 55.796
$ws3.Cells.Item(36,4).Value = -1.45
$ws3.Cells.Item(36,5).Value = 4.08
$ws3.Cells.Item(36,6).Value = "Grub"

$ws3.Cells.Item(37,1).Value = 55.788
$ws3.Cells.Item(37,4).Value = -1.57
$ws3.Cells.Item(37,5).Value = 3.98
$ws3.Cells.Item(37,6).Value = "Grub"

$ws3.Cells.Item(38,1).Value = 55.739
$ws3.Cells.Item(38,5).Value = 3.84
$ws3.Cells.Item(38,6).Value = "Grub"

$ws3.Cells.Item(39,1).Value = 55.715
$ws3.Cells.Item(39,2).Value = 15.39
$ws3.Cells.Item(39,3).Value = 0.095
$ws3.Cells.Item(39,4).Value = -1.65
$ws3.Cells.Item(39,5).Value = 4.3
$ws3.Cells.Item(39,6).Value = "Grub"

$ws3.Cells.Item(40,1).Value = 55.71
$ws3.Cells.Item(40,2).Value = 15.4
$ws3.Cells.Item(40,3).Value = 0.095
$ws3.Cells.Item(40,4).Value = -1.67
$ws3.Cells.Item(40,5).Value = 3.65
$ws3.Cells.Item(40,6).Value = "Grub"

$ws3.Cells.Item(41,1).Value = 55.45
$ws3.Cells.Item(41,2).Value = 15.22
$ws3.Cells.Item(41,3).Value = 0.265
$ws3.Cells.Item(41,4).Value = -1.6
$ws3.Cells.Item(41,5).Value = 3.62
$ws3.Cells.Item(41,6).Value = "Grub"

# --- final view/selection state ---
$ws2.Range("A1:F49").Select()
$ws3.Select()
$ws3.Range("A19:XFD19").Select()
